$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the data rows (rows 2-10 keep their row numbers, values refreshed)
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = "15-JAN-26"
$ws.Range("B2").Value = "SM-452"
$ws.Range("C2").Value = "Air Arabia Egypt E5-394"
$ws.Range("D2").Value = 646
$ws.Range("E2").Value = 586
$ws.Range("F2").Value = 60

# Row 3
$ws.Range("A3").Value = "15-JAN-26"
$ws.Range("B3").Value = "SM-452"
$ws.Range("C3").Value = "Air Arabia Egypt E5-326"
$ws.Range("D3").Value = 702
$ws.Range("E3").Value = 586
$ws.Range("F3").Value = 116
$ws.Range("G3").Value = 40
$ws.Range("I3").Value = -10

# Row 4
$ws.Range("A4").Value = "19-JAN-26"
$ws.Range("C4").Value = "Air Arabia Egypt E5-324"
$ws.Range("D4").Value = 608
$ws.Range("E4").Value = 526
$ws.Range("F4").Value = 82
$ws.Range("G4").Value = 40
$ws.Range("I4").Value = -10

# Row 5
$ws.Range("D5").Value = 695
$ws.Range("E5").Value = 526
$ws.Range("F5").Value = 169

# Row 6
$ws.Range("A6").Value = "24-JAN-26"
$ws.Range("C6").Value = "Air Arabia Egypt E5-364"
$ws.Range("D6").Value = 633
$ws.Range("E6").Value = 781
$ws.Range("F6").Value = -148
$ws.Range("G6").Value = 40
$ws.Range("I6").Value = -10
$ws.Range("J6").Value = "MEDIUM THREAT - MONITOR"

# Row 7
$ws.Range("A7").Value = "24-JAN-26"
$ws.Range("C7").Value = "Air Arabia Egypt E5-326"
$ws.Range("D7").Value = 670
$ws.Range("E7").Value = 781
$ws.Range("F7").Value = -111
$ws.Range("G7").Value = 40
$ws.Range("I7").Value = -10
$ws.Range("J7").Value = "MEDIUM THREAT - MONITOR"

# Row 8
$ws.Range("A8").Value = "24-JAN-26"
$ws.Range("C8").Value = "flyadeal F3-773"
$ws.Range("D8").Value = 679
$ws.Range("E8").Value = 781
$ws.Range("F8").Value = -102
$ws.Range("G8").Value = 32
$ws.Range("I8").Value = -2

# Row 9
$ws.Range("A9").Value = "24-JAN-26"
$ws.Range("C9").Value = "Saudia SV-331"
$ws.Range("D9").Value = 774
$ws.Range("E9").Value = 781
$ws.Range("F9").Value = -7
$ws.Range("G9").Value = 46
$ws.Range("I9").Value = -16
$ws.Range("J9").Value = "MEDIUM THREAT - MONITOR"

# Row 10
$ws.Range("A10").Value = "24-JAN-26"
$ws.Range("C10").Value = "Saudia SV-411"
$ws.Range("D10").Value = 1006
$ws.Range("E10").Value = 781
$ws.Range("F10").Value = 225
$ws.Range("G10").Value = 46
$ws.Range("I10").Value = -16

# ---------------------------------------------------------------------------
# 2. Remove the now-obsolete trailing rows (old rows 11-16)
# ---------------------------------------------------------------------------
$ws.Rows("11:16").Delete()

# ---------------------------------------------------------------------------
# 3. Widen column J (10th column) so the longer "MEDIUM THREAT - MONITOR"
#    label fits (stored width 25 -> ColumnWidth property is offset by ~0.83)
# ---------------------------------------------------------------------------
$ws.Columns(10).ColumnWidth = 24.17

# ---------------------------------------------------------------------------
# 4. Give the new "MEDIUM THREAT - MONITOR" cells their own highlight style
#    (new fill colour FFF3CD) while keeping the same bold font / border /
#    centered alignment as the existing LOW THREAT style.
# ---------------------------------------------------------------------------
$ws.Range("J6").Interior.Color = 13497343
$ws.Range("J7").Interior.Color = 13497343
$ws.Range("J9").Interior.Color = 13497343
